# Applies the following changes to Playwright_framework.pptx:
#  1. Slide 15: remove the empty tooltip="" from the "Playwright Trace Viewer" hyperlink.
#  2. Slide 2 (Agenda): "Custom Execution - Browserstack & LambdaTest" -> "Custom Execution - Browserstack"
#  3. Slide 8 title: "Custom Execution - Browserstack & LambdaTest" -> "Custom Execution - Browserstack"
#  4. Slide 8 body: append " in .env file" to "BROWSERSTACK_ACCESS_KEY" and remove the whole
#     "Execution in LambdaTest" section (the trailing empty paragraph + the LambdaTest header
#     paragraph + the LambdaTest env-vars paragraph).

$p = $ppt.ActivePresentation

# --- 1. Slide 15: drop the empty tooltip on the Playwright Trace Viewer hyperlink ---
$s15 = $p.Slides.Item(15)
$shp15 = $s15.Shapes.Item(2)
$tr15 = $shp15.TextFrame.TextRange
$text15 = $tr15.Text
$linkText = "Playwright Trace Viewer"
$idx15 = $text15.IndexOf($linkText)
$linkRange = $tr15.Characters($idx15 + 1, $linkText.Length)
$clickSetting = $linkRange.ActionSettings.Item(1)
$clickSetting.Hyperlink.ScreenTip = $null

# --- 2. Slide 2 (Agenda): drop "& LambdaTest" from the agenda bullet ---
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(2)
$tr2 = $shp2.TextFrame.TextRange
$oldAgenda = "Custom Execution - Browserstack & LambdaTest"
$newAgenda = "Custom Execution - Browserstack"
$text2 = $tr2.Text
$idx2 = $text2.IndexOf($oldAgenda)
$agendaRange = $tr2.Characters($idx2 + 1, $oldAgenda.Length)
$agendaRange.Text = $newAgenda

# --- 3. Slide 8 title: drop "& LambdaTest" ---
$s8 = $p.Slides.Item(8)
$titleShp8 = $s8.Shapes.Item(1)
$titleTr8 = $titleShp8.TextFrame.TextRange
$oldTitle = "Custom Execution - Browserstack & LambdaTest"
$newTitle = "Custom Execution - Browserstack"
$titleText = $titleTr8.Text
$idxTitle = $titleText.IndexOf($oldTitle)
$titleRange = $titleTr8.Characters($idxTitle + 1, $oldTitle.Length)
$titleRange.Text = $newTitle

# --- 4. Slide 8 body: remove the LambdaTest section, extend the Browserstack note ---
$bodyShp8 = $s8.Shapes.Item(2)
$bodyTr8 = $bodyShp8.TextFrame.TextRange
$keyMarker = "BROWSERSTACK_ACCESS_KEY"

$bodyText = $bodyTr8.Text
$keyIdx = $bodyText.IndexOf($keyMarker)
$deleteStart = $keyIdx + $keyMarker.Length
$deleteLen = $bodyText.Length - $deleteStart + 1
$tailRange = $bodyTr8.Characters($deleteStart + 1, $deleteLen)
$tailRange.Delete()

$bodyText2 = $bodyTr8.Text
$keyIdx2 = $bodyText2.IndexOf($keyMarker)
$keyRange = $bodyTr8.Characters($keyIdx2 + 1, $keyMarker.Length)
$keyRange.Text = "$keyMarker in .env file"
